$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing table from rows 13:16 up to rows 1:4 ---
$ws.Range("A13:C16").Copy()
$ws.Range("A1").PasteSpecial()
$ws.Range("A13:C16").Clear()

# --- Add the new "Section 2" row of data ---
$ws.Range("B4").Value = "Section 2: Creating with API requests"
$ws.Range("C4").Value = 3

# --- Update the selection / view state ---
$ws.Range("F11").Select()
